$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text representation
# (Excel would otherwise auto-convert numeric-looking strings like "1.002"
# into real numbers, which would lose trailing zeros / change formatting).
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.696.99"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.877.90"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "330.63"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "0.4724"
$ws.Range("E7").Value = "  +5.06%  "
$ws.Range("D8").Value = "0.3971"
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("D9").Value = "47.80"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "0.08051"
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").Value = "1.025"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "21.90"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "1.852.35"
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "5.970"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "7.181"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "87.17"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "0.00001047"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "0.06614"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Value = "17.33"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "27.709.57"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("D23").Value = "5.514"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "2.299"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "2.098.43"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").Value = "156.39"
$ws.Range("E27").Value = "  +3.45%  "
$ws.Range("D28").Value = "20.29"
$ws.Range("E28").Value = "  +4.11%  "
$ws.Range("D29").Value = "2.097"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("D30").Value = "5.599"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").Value = "122.58"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").Value = "0.9727"
$ws.Range("E32").Value = "  +4.64%  "
$ws.Range("D33").Value = "0.09564"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("D34").Value = "1.451"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").Value = "3.626"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "5.328"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("D37").Value = "0.06120"
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "1.236"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "8.157"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").Value = "0.1907"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("D44").Value = "10.24"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").Value = "0.5723"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").Value = "1.252"
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "3.401"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("D49").Value = "1.936"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "0.00000000317"
$ws.Range("E50").Value = "  +10.36%  "
$ws.Range("D51").Value = "0.06818"
$ws.Range("E51").Value = "  -0.59%  "

# Restore the default (General) style on the range so no stray cell-level
# style references are introduced by the temporary text number format.
$priceVolRange.Style = "Normal"

